# ---------------------------------------------------------------------------
# Refresh the Price (col D) and Volume(1h) (col E) figures on the crypto symbol
# list, 8-2-2023 12:00 GitHub Actions run.
#
# The source sheet stores every data cell as literal text (coinranking.com
# scrape dumped straight to XML as inline strings) - "332.03", "0.67%", etc are
# not numbers/percentages to Excel, just strings that happen to look numeric.
# Writing a numeric-looking string into a General-formatted cell makes Excel
# silently reinterpret it as a Number (and percentages as fractions styled with
# "0%"), which would change the stored cell type. Forcing the cell to the Text
# number format ("@") before the write keeps it a plain string, matching the
# source data untouched by this update.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    'D2' = '332.11'
    'E2' = '0.60%'
    'D3' = '45.43'
    'E3' = '2.44%'
    'D4' = '5.519'
    'E4' = '0.41%'
    'D5' = '0.08459'
    'E5' = '5.39%'
    'D6' = '2.168'
    'E6' = '5.03%'
    'D7' = '0.9901'
    'E7' = '3.75%'
    'D8' = '2.554'
    'E8' = '-0.07%'
    'D9' = '0.1189'
    'E9' = '4.18%'
    'D10' = '0.1937'
    'E10' = '2.56%'
    'D11' = '9.505'
    'E11' = '-6.65%'
    'D12' = '0.09862'
    'E12' = '-0.65%'
    'D13' = '0.04680'
    'E13' = '-4.03%'
    'E14' = '-0.37%'
    'D15' = '0.001290'
    'E15' = '1.67%'
    'D16' = '0.005887'
    'E16' = '2.00%'
    'D17' = '3.389'
    'E17' = '-0.05%'
    'D18' = '4.425'
    'E18' = '0.40%'
    'D19' = '0.3338'
    'E19' = '-1.67%'
    'D20' = '0.1355'
    'E20' = '-1.98%'
    'D21' = '0.2544'
    'E21' = '-1.45%'
    'D22' = '0.04144'
    'E22' = '1.42%'
    'D23' = '0.001296'
    'E23' = '-0.39%'
    'D24' = '0.004568'
    'E24' = '4.93%'
    'D25' = '0.0001304'
    'E25' = '8.74%'
    'E26' = '0.06%'
    'D38' = '0.02726'
    'E38' = '4.87%'
    'D39' = '0.05800'
    'E39' = '-0.12%'
    'D40' = '0.007888'
    'E40' = '4.28%'
    'D41' = '0.1436'
    'E41' = '2.26%'
    'D42' = '0.007515'
    'E42' = '2.60%'
    'D43' = '0.002026'
    'E43' = '0.62%'
    'D44' = '0.008957'
    'E44' = '9.29%'
    'D45' = '0.3546'
    'D46' = '0.00007115'
    'E46' = '1.36%'
    'D47' = '0.00000000752'
    'E47' = '0.26%'
    'E48' = '0.32%'
    'D49' = '0.003538'
    'E49' = '0.26%'
    'D50' = '0.003007'
    'E50' = '-14.15%'
    'D51' = '0.00002105'
    'E51' = '0.26%'
}

foreach ($cellRef in $updates.Keys) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$cellRef]
}
